# Fix shark double counts: update computed percentage values in rows
# 10, 11, 12, 15, 16, and 20 for columns B through G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    10 = @{ B = 95.38283744105355; C = 6.341402055340052; D = 39.51851520369351; E = 54.14008274096643; F = 45.85991725903357; G = 54.14008274096643 }
    11 = @{ B = 95.35980609945506; C = 25.99140387354132; D = 44.04745699154409; E = 29.9611391349146;  F = 70.03886086508541; G = 29.9611391349146  }
    12 = @{ B = 98.4871602932466;  C = 14.93902192521486; D = 67.50435457160448; E = 17.55662350318067; F = 82.44337649681934; G = 17.55662350318067 }
    15 = @{ B = 98.26681293112539; C = 22.93056935485755; D = 26.57459035111507; E = 50.49484029402739; F = 49.50515970597262; G = 50.49484029402739 }
    16 = @{ B = 87.65426515699606; C = 43.27108330943184; D = 38.29205176940727; E = 18.43686492116091; F = 81.56313507883911; G = 18.43686492116091 }
    20 = @{ B = 88.70339407047719; C = 26.52088558783935; D = 48.27323632967825; E = 25.2058780824824;  F = 74.7941219175176;  G = 25.2058780824824  }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}
